$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2778.902526399997
$ws.Range("E2").Value = 290927.2506141524
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312751
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 71977.22211760026
$ws.Range("O2").Value = 68708.80120585456

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52443.38413002542
$ws.Range("E2").Value = 269427.8185625125
$ws.Range("I2").Value = 226515.4616545431
$ws.Range("L2").Value = 216244.0719986508
$ws.Range("M2").Value = 105684.9748647551
$ws.Range("N2").Value = 35977.50527378646
$ws.Range("O2").Value = 25179.37288612182

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22664.8797166071
$ws.Range("B2").Value = 15375.40221492003
$ws.Range("E2").Value = 110456.9732286501
$ws.Range("I2").Value = 162836.6223031954
$ws.Range("M2").Value = 58634.61533225987
$ws.Range("N2").Value = 49759.46314991158
$ws.Range("O2").Value = 58375.08308716356
